# Edit LOT2064.docx: split the three long run-on paragraphs ("Programa" PT,
# "Programa" EN/italic, and "Bibliografia") into per-item runs separated by
# manual line breaks (<w:br/>), matching the target revision.

$d = $word.ActiveDocument

function Split-ParagraphIntoBreaks {
    param(
        [string]$AnchorText,
        [string]$PackageXml
    )

    $rng = $d.Content
    $found = $rng.Find.Execute($AnchorText, $true, $false, $false, $false, $false, $true, 1, $false)
    if (-not $found) {
        throw "Anchor text not found: $AnchorText"
    }

    $para = $rng.Paragraphs(1)
    $prange = $para.Range
    # Exclude the trailing paragraph-mark character from the range so the
    # paragraph itself is preserved (only its run content is replaced).
    $prange.MoveEnd(1, -1)
    $prange.Delete()
    $prange.InsertXML($PackageXml)
}

# --- Paragraph 1: "Programa" section, Portuguese -------------------------
$xml1 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>1)EMPREENDEDORISMO CIENTÍFICO E INOVAÇÃO: conceitos e definições; o perfil empreendedor, modelo de negócios; empreendedorismo de base científica e tecnológica; ecossistemas de inovação; start-ups, aceleradoras e corporate vendure; estratégias e desafios para abertura de startups de base científica; interação universidade-empresa; fases do desenvolvimento tecnológico e industrial;</w:t><w:br/><w:t xml:space="preserve">2)ADMINSITRAÇÃO DE EMPRESAS: conceitos fundamentais de administração de empresas e as formas de configuração de sua estrutura organizacional </w:t><w:br/><w:t>3)PROPRIEDADE INTELECTURAL: formas de proteção; busca de anterioridade; como solicitar patente; transferência de tecnologia</w:t><w:br/><w:t>4)PLANEJAMENTO MERCADOLÓGICO: Pesquisa e análise de mercado. Estratégias e ações para atingir o mercado. Análise Concorrência. Como criar negócios de alto crescimento;</w:t><w:br/><w:t>5)PLANEJAMENTO FINANCEIRO: Estimativa de Receitas e Custos da operação; Custos para abertura da startup; Formas para financiar o início do seu negócio; Fundos para financiamento do desenvolvimento de tecnologias de base científica</w:t><w:br/><w:t>6)FORMAÇÃO E GESTÃO DE EQUIPES: Construção e gestão de equipes de sucesso;</w:t><w:br/><w:t>7)VENDENDO SEU PRODUTO/SERVIÇO: Principais formas para inovar no seu negócio: Processo, Produto/Serviço, Posicionamento e Modelo de Negócio. Como elaborar uma estratégia futura. Estratégia de vendas. Pitch de vendas. Ciclo de vendas; Relacionamento com cliente;</w:t><w:br/><w:t>8)PLANO DE NEGÓCIOS: Plano de Negócios. Fluxo de Caixa. Lucratividade. Caixa: projeção de vendas. Aspectos legais para a abertura do seu negócio; métricas para startups</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Split-ParagraphIntoBreaks "1)EMPREENDEDORISMO CIENTÍFICO E INOVAÇÃO" $xml1

# --- Paragraph 2: "Programa" section, English/italic ----------------------
$xml2 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t>1.ENTREPRENEURSHIP AND INNOVATION: concepts and definitions; the entrepreneurial profile, business model; technology-based entrepreneurship; innovation ecosystems; start-ups, accelerators and corporate vendure</w:t><w:br/><w:t>2.INTELLECTURAL PROPERTY: forms of protection; how to apply for patent</w:t><w:br/><w:t>3.MARKET PLANNING: Market research and analysis; strategies and actions to reach the market; Competition analysis; gow to create high growth businesses;</w:t><w:br/><w:t>4.FINANCIAL PLANNING: estimated revenue and costs of the operation; costs for create a startup; ways to finance the start of your business;</w:t><w:br/><w:t>5.TEAM TRAINING; Building successful teams;</w:t><w:br/><w:t>6.SELLING YOUR PRODUCT / SERVICE: Forms to innovate in your business: Process, Product / Service, Positioning, and Business Model. How to elaborate a future strategy. sales strategy. Pitch of sales. Sales cycle; customer relationship;</w:t><w:br/><w:t>7.BUSINESS PLAN: business plan; Cash flow; profitability; cashier: sales projection; legal aspects for opening your business; metrics for startups.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Split-ParagraphIntoBreaks "1.ENTREPRENEURSHIP AND INNOVATION" $xml2

# --- Paragraph 3: "Bibliografia" section -----------------------------------
$xml3 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">•BLANK, Steve Gary. Do Sonho a realização em 4 passos: Estratégias para a criação de empresas de sucesso. Editora Evora. 3ª edição, 2008 </w:t><w:br/><w:t xml:space="preserve">•BLANK, Steve; DORF, Bob. STARTUP: Manual do Empreendedorismo. O guia passo a passo para construir uma grande empresa. Alta Books Editora. 1ª edição, 2014. </w:t><w:br/><w:t xml:space="preserve">•CECCONELO, Antonio; AJZENTAL, Alberto. A construção do plano de negócios. Ed. Saraiva, 1ª edição, 2008. </w:t><w:br/><w:t xml:space="preserve">•CHIAVENATO, Idalberto. Empreendedorismo – dando asas ao espírito empreendedor. Ed. Saraiva, 3ª edição, 2008. </w:t><w:br/><w:t xml:space="preserve">•DORNELAS, Jose. Empreendedorismo: transformando ideias em negócios. Editora Campus. 1ª edição, 2001 </w:t><w:br/><w:t xml:space="preserve">•DORNELAS, Jose. Empreendedorismo na prática. LTC. 3ª edição, 2015 </w:t><w:br/><w:t xml:space="preserve">•DORNELAS, Jose Carlos Assis. Empreendedorismo na prática – mitos e verdades do empreendedor de sucesso. Elsevier/Campus: Rio de Janeiro, 2007. </w:t><w:br/><w:t xml:space="preserve">•FILION, L. J.; Visão e Relações: Elementos para um Metamodelo da Atividade Empreendedora. International Small Business Journal, 1991. Tradução de Costa, S.R. </w:t><w:br/><w:t xml:space="preserve">•FILION, L. J.; O planejamento do seu Sistema de Aprendizagem Empresarial: Identifique uma Visão e Avalie o seu Sistema de Relações. Revista de Administração de Empresas, FGV, São Paulo, jul/set. 1991, pag. 31(3): 63:71. </w:t><w:br/><w:t xml:space="preserve">•HASHIMOTO, Marcos. Espírito empreendedor nas organizações – aumentando a competitividade através do intraempreendedorismo. São Paulo: Saraiva, 2006. </w:t><w:br/><w:t xml:space="preserve">•HISRICH, Robert; PETERS, Michael. Empreendedorismo. 5.ed. - Porto Alegre: Bookman, 2004. </w:t><w:br/><w:t xml:space="preserve">•OSTERWALDER, Alexander. Inovação Em Modelos de Negócios – Business Model Generation. Editora Alta Books, 2011 </w:t><w:br/><w:t xml:space="preserve">•PINCHOT, Gifford; PELLMAN, Ron. Intraempreendedorismo na prática: um guia de inovação. Campus: 2004 </w:t><w:br/><w:t xml:space="preserve">•RIES, Eric. A startup enxuta. Leya Editora. 1ª edição, 2011 </w:t><w:br/><w:t xml:space="preserve">•SANTOS. S.A. e CUNHA, N.C.V (orgs.). Empresas de Base Tecnológica: Conceitos, instrumentos e recursos. Unicorpore, 2005 </w:t><w:br/><w:t xml:space="preserve">•THIEL, Peter. De Zero a UM: O que aprender sobre empreendedorismo com Vale do Silício. Objetiva. 1ª edição, 2014 </w:t><w:br/><w:t>•TIMMONS; Jeffry; DORNELAS, José. SPINELLI, Stephen. A criação de novos negócios – empreendedorismo para o século 21. Editora Campus. 2010.</w:t><w:br/><w:t>•Critical Success Factors for Technology-based Startups, Jose Santisteban and David Mauricio; Int. Journal of Enterpreneurship and Small Business, vol. 42, no. 4, Inderscience, 2021.</w:t><w:br/><w:t>•SILVA, M. M. L. Administração para Estudantes e Profissionais de Áreas Técnicas. São Paulo: Brasport, 2018.</w:t><w:br/><w:t>•GUERRINI, F. M.; ESCRIÇÃO FILHO, E.; ROSIM, D.: Administração Para Engenheiros. Rio de Janeiro: Campus, 2016.</w:t><w:br/><w:t>•CHIAVENATO, I: Introdução À Teoria Geral da Administração. 9 ed. São Paulo: Manole, 2014.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
Split-ParagraphIntoBreaks "•BLANK, Steve Gary" $xml3

Write-Output "Done."
